# Apply the edit described by the diff:
#  - Populate column P (new values) and column S (new values) for rows 8-27
#  - Update the sheet view: topLeftCell="F14", selection activeCell="P28" sqref="P28"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    8  = @{ P = 34; S = 53 }
    9  = @{ P = 48; S = 60 }
    10 = @{ P = 64; S = 66 }
    11 = @{ P = 70; S = 65 }
    12 = @{ P = 74; S = 71 }
    13 = @{ P = 73; S = 72 }
    14 = @{ P = 78; S = 75 }
    15 = @{ P = 78; S = 72 }
    16 = @{ P = 79; S = 74 }
    17 = @{ P = 81; S = 76 }
    18 = @{ P = 80; S = 77 }
    19 = @{ P = 80; S = 77 }
    20 = @{ P = 78; S = 77 }
    21 = @{ P = 79; S = 77 }
    22 = @{ P = 80; S = 77 }
    23 = @{ P = 81; S = 78 }
    24 = @{ P = 80; S = 77 }
    25 = @{ P = 81; S = 78 }
    26 = @{ P = 81; S = 77 }
    27 = @{ P = 81; S = 78 }
}

foreach ($row in $values.Keys) {
    $ws.Range("P$row").Value = $values[$row].P
    $ws.Range("S$row").Value = $values[$row].S
}

# Update the active window view/selection
$ws.Range("P28").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 6
